$wb = $excel.ActiveWorkbook

# ---- Sheet "M 2020" ----
$ws2020 = $wb.Worksheets.Item("M 2020")

# Row 18
$ws2020.Range("P18").Value = "WI/K/1, WI/K/2, WI/K/3, WI/V/2   "
$ws2020.Range("O18").Value = "ja"

# Row 19
$ws2020.Range("P19").Value = "WI/K/1, WI/K/2, WI/K/7, WI/V/3"
$ws2020.Range("O19").Value = "ja"

# Row 20
$ws2020.Range("G20").Value = 2
$ws2020.Range("H20").Value = "Hoofdstuk 4: Machtsverbanden, Hoofdstuk 5: Rekenen, Hoofdstuk 6: Goniometrie"
$ws2020.Range("J20").Value = "tt"
$ws2020.Range("L20").Value = 100
$ws2020.Range("N20").Value = 3
$ws2020.Range("O20").Value = "ja"
$ws2020.Range("P20").Value = "WI/K/1, WI/K/2, WI/K/3, WI/V/2, Rekenen"

# Row 21
$ws2020.Range("G21").Value = 3
$ws2020.Range("H21").Value = "Hoofdstuk 7: Exponentiële formules, Hoofdstuk 8: Ruimtemeetkunde"
$ws2020.Range("N21").Value = 2
$ws2020.Range("O21").Value = "ja"
$ws2020.Range("P21").Value = "WI/K/1, WI/K/2, WI/K/3, WI/V/2"

# Row 22 - cleared entirely (back to "kies..." placeholders)
$ws2020.Range("D22").ClearContents()
$ws2020.Range("G22").Value = "kies…"
$ws2020.Range("H22").ClearContents()
$ws2020.Range("J22").Value = "kies…"
$ws2020.Range("L22").ClearContents()
$ws2020.Range("M22").Value = "kies…"
$ws2020.Range("N22").ClearContents()
$ws2020.Range("O22").Value = "kies…"
$ws2020.Range("P22").ClearContents()

# ---- Sheet "M 2019" ----
$ws2019 = $wb.Worksheets.Item("M 2019")
$ws2019.Range("O18").Value = "ja"
$ws2019.Range("O19").Value = "ja"
# H20 text was stored with mangled (mojibake) encoding; rewrite with the
# correct UTF-8 text (content is otherwise unchanged).
$ws2019.Range("H20").Value = "opdracht over oriëntatie op leren en werken bij een zelfgekozen sector"
$ws2019.Range("O21").Value = "ja"
# H22 text had the same mojibake issue for the "ë"; correct it here too.
$ws2019.Range("H22").Value = "Hoofdstuk 7: Exponentiële formules, Hoofdstuk 8: Ruimtemeetkunde"
$ws2019.Range("O22").Value = "ja"
